$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows above the current row 2, pushing all existing data rows
# (originally rows 2-54) down to rows 8-60.
$ws.Range("A2:T7").EntireRow.Insert()

# The freshly inserted rows picked up formatting from the header row above;
# reset that, then restore the date number format on column D like the rest
# of the data rows.
$ws.Range("A2:T7").ClearFormats()
$ws.Range("D2:D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Columns that are identical across all six new rows.
$ws.Range("A2:A7").Value = 9
$ws.Range("B2:B7").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C2:C7").Value = "Metropolitana"
$ws.Range("D2:D7").Value = 44552
$ws.Range("E2:E7").Value = 13
$ws.Range("F2:F7").Value = "Fruta"
$ws.Range("G2:G7").Value = 100103
$ws.Range("H2:H7").Value = "Frutos de hueso (carozo)"
$ws.Range("I2:I7").Value = 100103003
$ws.Range("J2:J7").Value = "Damasco"
$ws.Range("R2:R7").Value = "Provincia de San Felipe de Aconcagua"

# Row 2: Dina / Especial
$ws.Range("K2").Value = "Dina"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 10000
$ws.Range("P2").Value = 10000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("S2").Value = 1000
$ws.Range("T2").Value = 10

# Row 3: Dina / Primera
$ws.Range("K3").Value = "Dina"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 220
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 800
$ws.Range("T3").Value = 10

# Row 4: Dina / Segunda
$ws.Range("K4").Value = "Dina"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 350
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("S4").Value = 600
$ws.Range("T4").Value = 10

# Row 5: Modesto / Especial
$ws.Range("K5").Value = "Modesto"
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 280
$ws.Range("N5").Value = 21600
$ws.Range("O5").Value = 21600
$ws.Range("P5").Value = 21600
$ws.Range("Q5").Value = "$/caja 18 kilos granel"
$ws.Range("S5").Value = 1200
$ws.Range("T5").Value = 18

# Row 6: Modesto / Primera
$ws.Range("K6").Value = "Modesto"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 310
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("Q6").Value = "$/caja 18 kilos granel"
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 18

# Row 7: Modesto / Segunda
$ws.Range("K7").Value = "Modesto"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 350
$ws.Range("N7").Value = 14400
$ws.Range("O7").Value = 14400
$ws.Range("P7").Value = 14400
$ws.Range("Q7").Value = "$/caja 18 kilos granel"
$ws.Range("S7").Value = 800
$ws.Range("T7").Value = 18
